$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 555.325
$ws.Range("J17").Value = 567.91895
$ws.Range("L17").Value = 1703.75685
$ws.Range("N17").Value = -2039.75685
$ws.Range("H64").Value = 4978.4
$ws.Range("J64").Value = 4297
$ws.Range("L64").Value = 4297
$ws.Range("N64").Value = -4793
$ws.Range("H67").Value = 4978.4
$ws.Range("J67").Value = 4297
$ws.Range("L67").Value = 4297
$ws.Range("N67").Value = -6013
$ws.Range("H74").Value = 143722.62
$ws.Range("I74").Value = 143722.62
$ws.Range("K74").Value = 143722.62
$ws.Range("M74").Value = -142786.62
$ws.Range("H77").Value = 143722.62
$ws.Range("I77").Value = 143722.62
$ws.Range("K77").Value = 718613.1
$ws.Range("M77").Value = -713933.1
$ws.Range("H106").Value = 13445.363
$ws.Range("I106").Value = 12666.5
$ws.Range("K106").Value = 12666.5
$ws.Range("M106").Value = -12035.5
$ws.Range("N110").Value = -108179
$ws.Range("H110").Value = 99999
$ws.Range("J110").Value = 99999
$ws.Range("L110").Value = 99999
$ws.Range("H129").Value = 4307.1055
$ws.Range("I129").Value = 996.0909
$ws.Range("K129").Value = 2988.2727
$ws.Range("M129").Value = 2011.7273
$ws.Range("H131").Value = 4565634
$ws.Range("I131").Value = 16966.334
$ws.Range("K131").Value = 50899.00199999999
$ws.Range("M131").Value = -45859.00199999999
$ws.Range("H138").Value = 3795.7307
$ws.Range("I138").Value = 1334.2916
$ws.Range("J138").Value = 33333
$ws.Range("K138").Value = 4002.8748
$ws.Range("L138").Value = 99999
$ws.Range("M138").Value = 1137.1252
$ws.Range("N138").Value = -110279
$ws.Range("H141").Value = 22736322
$ws.Range("I141").Value = 31253948
$ws.Range("J141").Value = 22656
$ws.Range("K141").Value = 93761844
$ws.Range("L141").Value = 67968
$ws.Range("M141").Value = -93756664
$ws.Range("N141").Value = -78328
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3013.606
$ws.Range("I32").Value = 2951.5312
$ws.Range("K32").Value = 2951.5312
$ws.Range("M32").Value = -2664.5312
$ws.Range("N34").Value = -90541
$ws.Range("H34").Value = 89999
$ws.Range("J34").Value = 89999
$ws.Range("L34").Value = 89999
$ws.Range("H124").Value = 30966.334
$ws.Range("J124").Value = 30966.334
$ws.Range("L124").Value = 30966.334
$ws.Range("N124").Value = -40786.334
$ws.Range("H132").Value = 3576403.5
$ws.Range("I132").Value = 5332.04
$ws.Range("J132").Value = 33335332
$ws.Range("K132").Value = 15996.12
$ws.Range("L132").Value = 100005996
$ws.Range("M132").Value = -13466.12
$ws.Range("N132").Value = -100011056
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1050.25
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 1167
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 1167
$ws.Range("M22").Value = -527
$ws.Range("N22").Value = -1513
$ws.Range("H86").Value = 34881.953
$ws.Range("I86").Value = 68985.5
$ws.Range("J86").Value = 3878.7273
$ws.Range("K86").Value = 68985.5
$ws.Range("L86").Value = 3878.7273
$ws.Range("M86").Value = -67862.5
$ws.Range("N86").Value = -6124.7273
$ws.Range("H89").Value = 34881.953
$ws.Range("I89").Value = 68985.5
$ws.Range("J89").Value = 3878.7273
$ws.Range("K89").Value = 344927.5
$ws.Range("L89").Value = 19393.6365
$ws.Range("M89").Value = -339311.5
$ws.Range("N89").Value = -30625.6365
$ws.Range("H105").Value = 850634.7
$ws.Range("I105").Value = 1041319.94
$ws.Range("J105").Value = 11619.8
$ws.Range("K105").Value = 1041319.94
$ws.Range("L105").Value = 11619.8
$ws.Range("M105").Value = -1039572.94
$ws.Range("N105").Value = -15113.8
$ws.Range("H134").Value = 6669654.5
$ws.Range("I134").Value = 2404.6
$ws.Range("K134").Value = 7213.799999999999
$ws.Range("M134").Value = -4678.799999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3192
$ws.Range("I122").Value = 3156.6
$ws.Range("K122").Value = 9469.799999999999
$ws.Range("M122").Value = -7019.799999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4579.033
$ws.Range("I131").Value = 3790
$ws.Range("J131").Value = 4819.174
$ws.Range("K131").Value = 11370
$ws.Range("L131").Value = 14457.522
$ws.Range("M131").Value = -6330
$ws.Range("N131").Value = -24537.522
$ws.Range("H140").Value = 6351.2856
$ws.Range("I140").Value = 1854.3334
$ws.Range("K140").Value = 5563.0002
$ws.Range("M140").Value = -383.0002000000004
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2716.3333
$ws.Range("J80").Value = 3088.2856
$ws.Range("L80").Value = 3088.2856
$ws.Range("N80").Value = -5084.2856
$ws.Range("H83").Value = 2716.3333
$ws.Range("J83").Value = 3088.2856
$ws.Range("L83").Value = 15441.428
$ws.Range("N83").Value = -25425.428
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 6662.6665
$ws.Range("I122").Value = 6662.6665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19987.9995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -17537.9995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2878.3845
$ws.Range("I22").Value = 2490.7144
$ws.Range("K22").Value = 2490.7144
$ws.Range("M22").Value = -2195.7144
$ws.Range("H27").Value = 2878.3845
$ws.Range("I27").Value = 2490.7144
$ws.Range("K27").Value = 2490.7144
$ws.Range("M27").Value = -2383.7144
$ws.Range("H55").Value = 1488.5454
$ws.Range("I55").Value = 1685
$ws.Range("J55").Value = 1396.8667
$ws.Range("K55").Value = 1685
$ws.Range("L55").Value = 1396.8667
$ws.Range("M55").Value = -1512
$ws.Range("N55").Value = -1742.8667
$ws.Range("H132").Value = 3983.9443
$ws.Range("I132").Value = 2555
$ws.Range("K132").Value = 7665
$ws.Range("M132").Value = -5135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 11725.2
$ws.Range("J96").Value = 14232
$ws.Range("L96").Value = 14232
$ws.Range("N96").Value = -16978
$ws.Range("H100").Value = 1782.8334
$ws.Range("I100").Value = 1799.4
$ws.Range("K100").Value = 3598.8
$ws.Range("M100").Value = -3057.8
$ws.Range("H122").Value = 3041.7144
$ws.Range("I122").Value = 2594
$ws.Range("K122").Value = 7782
$ws.Range("M122").Value = -5332
